$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in attendance/grade values of 5 for the given cells, preserving any
# existing cell formatting (fill color etc.) already applied to those cells.
$ws.Range("G12").Value = 5

$ws.Range("G20").Value = 5

$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 5

$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = 5

$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 5
$ws.Range("I26").Value = 5

$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 5
$ws.Range("I27").Value = 5

# Update the active cell/selection to match the recorded view state.
$ws.Range("G13").Select()
